$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.345.23"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "1.711.54"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("D4").Value = "'0.9957"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "'239.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.14%  "
$ws.Range("D6").Value = "'0.9965"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "'0.4885"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("D8").Value = "'0.2578"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.93%  "
$ws.Range("D9").Value = "'0.06171"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").Value = "1.713.92"
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("D11").Value = "'0.06945"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "'0.5977"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.99%  "
$ws.Range("D14").Value = "'4.456"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").Value = "'76.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "26.237.17"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").Value = "'0.9958"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "'0.000007070"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.29%  "
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("D21").Value = "1.933.11"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").Value = "'4.394"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.35%  "
$ws.Range("D23").Value = "'8.399"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.93%  "
$ws.Range("D24").Value = "'5.030"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.39%  "
$ws.Range("D25").Value = "'136.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("D26").Value = "'15.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.04%  "
$ws.Range("D27").Value = "'1.392"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'1.725"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("D29").Value = "'105.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.88%  "
$ws.Range("D30").Value = "'3.883"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("D31").Value = "'0.07920"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("D32").Value = "'3.599"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.56%  "
$ws.Range("D33").Value = "'0.04446"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("D34").Value = "'2.598"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").Value = "'0.9905"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.97%  "
$ws.Range("D36").Value = "'0.6173"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.31%  "
$ws.Range("D37").Value = "'0.9419"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.42%  "
$ws.Range("D38").Value = "'2.005"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").Value = "'0.9953"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").Value = "'0.01473"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.16%  "
$ws.Range("D42").Value = "'99.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("D43").Value = "'5.359"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("E44").Value = "  -3.49%  "
$ws.Range("D45").Value = "'6.796"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("D46").Value = "'0.1147"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("D47").Value = "'0.05343"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").Value = "'30.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").Value = "'7.721"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").Value = "'51.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("B51").Value = "TrueUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D51").Value = "'0.9988"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.33%  "
